$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.451.87"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.88%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.524.60"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.16%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "203.30"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "552.97"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.89%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.516.47"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.605"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.50%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.660"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "61.38"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +10.50%  "

$ws.Range("E12").Value = "  -6.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000275"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -5.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.85"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.085.46"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.510.80"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.31%  "

$ws.Range("E17").Value = "  -1.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.59"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "67.134.46"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.92"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.68%  "

$ws.Range("E21").Value = "  -4.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.96"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.31%  "

$ws.Range("E23").Value = "  -5.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.05"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -9.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.92"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.81"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.09"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.38%  "

$ws.Range("E28").Value = "  -4.37%  "

$ws.Range("E29").Value = "  -3.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.89"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.33"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -10.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "686.30"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.82"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.85%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "64.07"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.111"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.76"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -7.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.411"
$ws.Range("D37").ClearFormats()

$ws.Range("E38").Value = "  -0.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.07"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.095.22"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.00%  "

$ws.Range("E41").Value = "  -2.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0706"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -13.40%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.84"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +11.60%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.57"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -12.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.71"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +6.97%  "

$ws.Range("E47").Value = "  -4.53%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.128"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.49%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "136.84"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.56%  "

$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.33"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -6.48%  "

$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.95"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.45%  "
